$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("milestone")
Write-Host $ws.Name
